$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Bad Drivers summary (rows 3-5) ---
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 784
$ws.Range("D3").Value = 95.90000000000001

$ws.Range("D4").Value = 98.2

$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 853

# --- Good Drivers table (rows 13-21) ---
# New order/content for Adapter-Driver (A), Total Samples (B),
# Good Roaming Calculation % (D), Driver Vintage (E)

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B13").Value = 10661
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = 0

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B14").Value = 56018
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 0

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B15").Value = 34244
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = 0

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B16").Value = 442178
$ws.Range("D16").Value = 100
# Force text so the date-like string isn't auto-converted to a date serial,
# then restore the original (General) number format via a formats-only
# paste from an untouched style-4 cell, so the cell keeps style 4.
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2024-11-10"
$ws.Range("D13").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B17").Value = 14239
$ws.Range("D17").Value = 100
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2022-05-23"
$ws.Range("D13").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B18").Value = 265400
$ws.Range("D18").Value = 99.90000000000001
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2022-05-01"
$ws.Range("D13").Copy()
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B19").Value = 77849
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2021-08-18"
$ws.Range("D13").Copy()
$ws.Range("E19").PasteSpecial(-4122)

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B20").Value = 59673
$ws.Range("D20").Value = 100
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2020-08-05"
$ws.Range("D13").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B21").Value = 113652
$ws.Range("D21").Value = 100
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2019-12-14"
$ws.Range("D13").Copy()
$ws.Range("E21").PasteSpecial(-4122)
